$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (Source | C | FFR | LF) -----------------------------------
$ws.Range("B1").Value = "C"
$ws.Range("C1").Value = "FFR"
$ws.Range("D1").Value = "LF"
$ws.Range("A1").Value = "Source"

# --- Row labels (column A) -------------------------------------------------
$ws.Range("A2").Value = "C Lag"
$ws.Range("A3").Value = "FFR Lag"
$ws.Range("A4").Value = "LF Lag"

# --- Column B data (C Lag row) --------------------------------------------
$ws.Range("B2").Value = "-0.46***"
$ws.Range("B3").Value = "'-0.01"
$ws.Range("B4").Value = "0.04*"

# --- Column C data (FFR Lag row) -------------------------------------------
$ws.Range("C2").Value = "'3.79"
$ws.Range("C3").Value = "1.6***"
$ws.Range("C4").Value = "3.53*"

# --- Column D data (LF Lag row) --------------------------------------------
$ws.Range("D2").Value = "'-6.09"
$ws.Range("D3").Value = "0.5***"
$ws.Range("D4").Value = "0.54*"
